# Weekly fruit/vegetable data update: a new record is inserted as row 485
# (Feria Lagunitas de Puerto Montt - Betarraga), pushing all subsequent rows
# down by one (old row 485 -> 486, ..., old row 536 -> 537).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 485, shifting rows 485:536 down to 486:537.
$ws.Rows.Item(485).Insert()

# Populate the newly inserted row 485 with the new weekly record.
$ws.Cells.Item(485, 1).Value = 4
$ws.Cells.Item(485, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(485, 3).Value = "Los Lagos"
$ws.Cells.Item(485, 4).Value = 45194
$ws.Cells.Item(485, 5).Value = 10
$ws.Cells.Item(485, 6).Value = 100114014
$ws.Cells.Item(485, 7).Value = "Betarraga"
$ws.Cells.Item(485, 8).Value = "Sin especificar"
$ws.Cells.Item(485, 9).Value = "Primera"
$ws.Cells.Item(485, 10).Value = 500
$ws.Cells.Item(485, 11).Value = 1000
$ws.Cells.Item(485, 12).Value = 1000
$ws.Cells.Item(485, 13).Value = 1000
$ws.Cells.Item(485, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(485, 15).Value = "Región Metropolitana"
$ws.Cells.Item(485, 16).Value = 200
$ws.Cells.Item(485, 17).Value = 5
$ws.Cells.Item(485, 18).Value = "Hortaliza"
